$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# Force the "Week_Start_Date" / Summary "Value" cells we are about to rewrite to stay
# plain text (matching the source data, which stores these as inline-string dates and
# numeric-looking labels) instead of Excel auto-coercing "2025-02-02"-style strings
# into real dates.
$ws1.Range("B2:B17").NumberFormat = "@"
foreach ($r in @(2, 4, 5, 6, 8, 9, 10, 11, 12, 13, 14, 15)) {
    $ws2.Cells.Item($r, 2).NumberFormat = "@"
}

# --- Forecast Comparison: roll the 16-week forecast window forward by one week ---
# (new MyForecast / Amazon Mean / P70 / P80 / P90 figures from the refreshed model run)
$ws1.Range("B2").Value = "2025-02-02"
$ws1.Range("D2").Value = 196
$ws1.Range("E2").Value = 175
$ws1.Range("F2").Value = 204
$ws1.Range("G2").Value = 230
$ws1.Range("H2").Value = 270

$ws1.Range("B3").Value = "2025-02-09"
$ws1.Range("D3").Value = 169
$ws1.Range("E3").Value = 132
$ws1.Range("F3").Value = 158
$ws1.Range("G3").Value = 185
$ws1.Range("H3").Value = 228

$ws1.Range("B4").Value = "2025-02-16"
$ws1.Range("D4").Value = 123
$ws1.Range("E4").Value = 123
$ws1.Range("F4").Value = 148
$ws1.Range("G4").Value = 174
$ws1.Range("H4").Value = 214

$ws1.Range("B5").Value = "2025-02-23"
$ws1.Range("D5").Value = 92
$ws1.Range("E5").Value = 118
$ws1.Range("F5").Value = 142
$ws1.Range("G5").Value = 167
$ws1.Range("H5").Value = 205

$ws1.Range("B6").Value = "2025-03-02"
$ws1.Range("D6").Value = 95
$ws1.Range("E6").Value = 117
$ws1.Range("F6").Value = 141
$ws1.Range("G6").Value = 168
$ws1.Range("H6").Value = 209

$ws1.Range("B7").Value = "2025-03-09"
$ws1.Range("D7").Value = 122
$ws1.Range("E7").Value = 114
$ws1.Range("F7").Value = 137
$ws1.Range("G7").Value = 162
$ws1.Range("H7").Value = 201

$ws1.Range("B8").Value = "2025-03-16"
$ws1.Range("D8").Value = 140
$ws1.Range("E8").Value = 108
$ws1.Range("F8").Value = 131
$ws1.Range("G8").Value = 156
$ws1.Range("H8").Value = 195

$ws1.Range("B9").Value = "2025-03-23"
$ws1.Range("D9").Value = 135
$ws1.Range("E9").Value = 110
$ws1.Range("F9").Value = 133
$ws1.Range("G9").Value = 159
$ws1.Range("H9").Value = 200

$ws1.Range("B10").Value = "2025-03-30"
$ws1.Range("D10").Value = 116
$ws1.Range("E10").Value = 108
$ws1.Range("F10").Value = 131
$ws1.Range("G10").Value = 155
$ws1.Range("H10").Value = 193

$ws1.Range("B11").Value = "2025-04-06"
$ws1.Range("D11").Value = 123
$ws1.Range("E11").Value = 99
$ws1.Range("F11").Value = 120
$ws1.Range("G11").Value = 144
$ws1.Range("H11").Value = 181

$ws1.Range("B12").Value = "2025-04-13"
$ws1.Range("D12").Value = 130
$ws1.Range("E12").Value = 100
$ws1.Range("F12").Value = 121
$ws1.Range("G12").Value = 146
$ws1.Range("H12").Value = 186

$ws1.Range("B13").Value = "2025-04-20"
$ws1.Range("D13").Value = 125
$ws1.Range("E13").Value = 96
$ws1.Range("F13").Value = 117
$ws1.Range("G13").Value = 142
$ws1.Range("H13").Value = 181

$ws1.Range("B14").Value = "2025-04-27"
$ws1.Range("D14").Value = 125
$ws1.Range("E14").Value = 96
$ws1.Range("F14").Value = 117
$ws1.Range("G14").Value = 140
$ws1.Range("H14").Value = 176

$ws1.Range("B15").Value = "2025-05-04"
$ws1.Range("D15").Value = 116
$ws1.Range("E15").Value = 89
$ws1.Range("F15").Value = 109
$ws1.Range("G15").Value = 132
$ws1.Range("H15").Value = 169

$ws1.Range("B16").Value = "2025-05-11"
$ws1.Range("D16").Value = 118
$ws1.Range("E16").Value = 91
$ws1.Range("F16").Value = 111
$ws1.Range("G16").Value = 134
$ws1.Range("H16").Value = 171

$ws1.Range("B17").Value = "2025-05-18"
$ws1.Range("D17").Value = 94
$ws1.Range("E17").Value = 92
$ws1.Range("F17").Value = 112
$ws1.Range("G17").Value = 137
$ws1.Range("H17").Value = 175

# --- Summary: refreshed headline stats to match the updated forecast window ---
$ws2.Range("B2").Value = "2022-12-25 to 2025-01-26"
$ws2.Range("B4").Value = "208"
$ws2.Range("B5").Value = "88"
$ws2.Range("B6").Value = "81"
$ws2.Range("B8").Value = "9557 units"
$ws2.Range("B9").Value = "2019"
$ws2.Range("B10").Value = "1072"
$ws2.Range("B11").Value = "580"
$ws2.Range("B12").Value = "196"
$ws2.Range("B13").Value = "2025-02-02"
$ws2.Range("B14").Value = "92"
$ws2.Range("B15").Value = "2025-02-23"
